$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 130
$ws.Range("I33").Value = 119
$ws.Range("J33").Value = 185
$ws.Range("K33").Value = 119
$ws.Range("L33").Value = 185
$ws.Range("M33").Value = 110
$ws.Range("N33").Value = -643

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6916.25
$ws.Range("I62").Value = 5888.3335
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 5888.3335
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -5264.3335
$ws.Range("N62").Value = -11248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 6916.25
$ws.Range("I65").Value = 5888.3335
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 29441.6675
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -26321.6675
$ws.Range("N65").Value = -56240

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 800.8570999999999
$ws.Range("I132").Value = 825.5484
$ws.Range("K132").Value = 2476.6452
$ws.Range("M132").Value = 53.35480000000007

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1606.9286
$ws.Range("I135").Value = 1602.3636
$ws.Range("K135").Value = 14421.2724
$ws.Range("M135").Value = -11886.2724

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2076.238
$ws.Range("I137").Value = 904.2
$ws.Range("J137").Value = 3141.7273
$ws.Range("K137").Value = 2712.6
$ws.Range("L137").Value = 9425.1819
$ws.Range("M137").Value = -162.6000000000004
$ws.Range("N137").Value = -14525.1819

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7659.7617
$ws.Range("J138").Value = 8353.166999999999
$ws.Range("L138").Value = 25059.501
$ws.Range("N138").Value = -35339.501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2624.087
$ws.Range("I32").Value = 2387.158
$ws.Range("J32").Value = 3749.5
$ws.Range("K32").Value = 2387.158
$ws.Range("L32").Value = 3749.5
$ws.Range("M32").Value = -2100.158
$ws.Range("N32").Value = -4323.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5828.143
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2806.524
$ws.Range("I132").Value = 2269.4
$ws.Range("K132").Value = 6808.200000000001
$ws.Range("M132").Value = -4278.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5828.143
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3396.5
$ws.Range("I134").Value = 3095.2144
$ws.Range("J134").Value = 5505.5
$ws.Range("K134").Value = 9285.643199999999
$ws.Range("L134").Value = 16516.5
$ws.Range("M134").Value = -6750.643199999999
$ws.Range("N134").Value = -21586.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 19000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3051.28
$ws.Range("I132").Value = 2489.9524
$ws.Range("J132").Value = 5998.25
$ws.Range("K132").Value = 7469.8572
$ws.Range("L132").Value = 17994.75
$ws.Range("M132").Value = -4939.8572
$ws.Range("N132").Value = -23054.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3374.52
$ws.Range("I134").Value = 3348.4583
$ws.Range("K134").Value = 10045.3749
$ws.Range("M134").Value = -7510.374899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 799
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 799
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = 2397
$ws.Range("N5").Value = -2621

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1601.5
$ws.Range("I34").Value = 200
$ws.Range("K34").Value = 600
$ws.Range("M34").Value = -516

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("N39").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 382
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 289999.5
$ws.Range("I128").Value = 289999.5
$ws.Range("K128").Value = 869998.5
$ws.Range("M128").Value = -865018.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 799
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 799
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("M135").Value = 7191
$ws.Range("N135").Value = -12261

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5833
$ws.Range("I126").Value = 4999.5
$ws.Range("K126").Value = 14998.5
$ws.Range("M126").Value = -12528.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2614.8462
$ws.Range("I132").Value = 2416.1667
$ws.Range("K132").Value = 7248.500100000001
$ws.Range("M132").Value = -4718.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 889090.6
$ws.Range("I43").Value = 180000
$ws.Range("J43").Value = 922856.9
$ws.Range("K43").Value = 180000
$ws.Range("L43").Value = 922856.9
$ws.Range("M43").Value = -179807
$ws.Range("N43").Value = -923242.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2169.2856
$ws.Range("I46").Value = 858
$ws.Range("K46").Value = 858
$ws.Range("M46").Value = -670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 124.75
$ws.Range("I55").Value = 99.666664
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 99.666664
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = 73.333336
$ws.Range("N55").Value = -546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2955.5715
$ws.Range("J61").Value = 3144.5
$ws.Range("L61").Value = 3144.5
$ws.Range("N61").Value = -3548.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2024
$ws.Range("I93").Value = 2346.25
$ws.Range("J93").Value = 735
$ws.Range("K93").Value = 2346.25
$ws.Range("L93").Value = 735
$ws.Range("M93").Value = -1098.25
$ws.Range("N93").Value = -3231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2955.5715
$ws.Range("J113").Value = 3144.5
$ws.Range("L113").Value = 3144.5
$ws.Range("N113").Value = -7484.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2105.7273
$ws.Range("I132").Value = 1798.8334
$ws.Range("J132").Value = 3486.75
$ws.Range("K132").Value = 5396.5002
$ws.Range("L132").Value = 10460.25
$ws.Range("M132").Value = -2866.5002
$ws.Range("N132").Value = -15520.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1039
$ws.Range("I126").Value = 1451
$ws.Range("J126").Value = 489.66666
$ws.Range("K126").Value = 4353
$ws.Range("L126").Value = 1468.99998
$ws.Range("M126").Value = -1883
$ws.Range("N126").Value = -6408.999980000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1285.68
$ws.Range("I132").Value = 855.2632
$ws.Range("K132").Value = 2565.7896
$ws.Range("M132").Value = -35.78960000000006

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2366.2856
$ws.Range("I136").Value = 1599.1666
$ws.Range("K136").Value = 4797.4998
$ws.Range("M136").Value = -2247.4998
